$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 527
$ws.Range("I18").Value = 405.375
$ws.Range("K18").Value = 405.375
$ws.Range("M18").Value = -121.375
$ws.Range("H31").Value = 131.75
$ws.Range("I31").Value = 131.75
$ws.Range("K31").Value = 395.25
$ws.Range("M31").Value = -165.25
$ws.Range("H51").Value = 53999.19
$ws.Range("I51").Value = 9227.714
$ws.Range("J51").Value = 76384.93
$ws.Range("K51").Value = 9227.714
$ws.Range("L51").Value = 76384.93
$ws.Range("M51").Value = -8743.714
$ws.Range("N51").Value = -77352.93
$ws.Range("H64").Value = 5000
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("M64").Value = -4752
$ws.Range("H67").Value = 5000
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("M67").Value = -4142
$ws.Range("H76").Value = 499
$ws.Range("J76").Value = 499
$ws.Range("L76").Value = 499
$ws.Range("N76").Value = -1129
$ws.Range("H79").Value = 499
$ws.Range("J79").Value = 499
$ws.Range("L79").Value = 499
$ws.Range("N79").Value = -2683
$ws.Range("H92").Value = 389.75
$ws.Range("I92").Value = 308.91666
$ws.Range("K92").Value = 308.91666
$ws.Range("M92").Value = 939.08334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3929
$ws.Range("I32").Value = 2535.7942
$ws.Range("K32").Value = 2535.7942
$ws.Range("M32").Value = -2248.7942
$ws.Range("H34").Value = 25000
$ws.Range("J34").Value = 25000
$ws.Range("L34").Value = 25000
$ws.Range("N34").Value = -25542
$ws.Range("H61").Value = 10627
$ws.Range("I61").Value = 11430.857
$ws.Range("K61").Value = 11430.857
$ws.Range("M61").Value = -11218.857
$ws.Range("H74").Value = 3563.5
$ws.Range("J74").Value = 3995.6667
$ws.Range("L74").Value = 3995.6667
$ws.Range("N74").Value = -5743.6667
$ws.Range("H77").Value = 3563.5
$ws.Range("J77").Value = 3995.6667
$ws.Range("L77").Value = 19978.3335
$ws.Range("N77").Value = -28714.3335
$ws.Range("H136").Value = 10627
$ws.Range("I136").Value = 11430.857
$ws.Range("K136").Value = 34292.571
$ws.Range("M136").Value = -31742.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1592.6154
$ws.Range("I20").Value = 1651.3334
$ws.Range("J20").Value = 888
$ws.Range("K20").Value = 1651.3334
$ws.Range("L20").Value = 888
$ws.Range("M20").Value = -1404.3334
$ws.Range("N20").Value = -1382
$ws.Range("H39").Value = 1850
$ws.Range("I39").Value = 1500
$ws.Range("J39").Value = 5000
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 5000
$ws.Range("M39").Value = -1111
$ws.Range("N39").Value = -5778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 4595.4165
$ws.Range("I3").Value = 350.66666
$ws.Range("J3").Value = 8840.167
$ws.Range("K3").Value = 350.66666
$ws.Range("L3").Value = 8840.167
$ws.Range("M3").Value = -237.66666
$ws.Range("N3").Value = -9066.167
$ws.Range("H29").Value = 29500
$ws.Range("I29").Value = 29000
$ws.Range("J29").Value = 30000
$ws.Range("K29").Value = 29000
$ws.Range("L29").Value = 30000
$ws.Range("M29").Value = -28707
$ws.Range("N29").Value = -30586
$ws.Range("H58").Value = 3273
$ws.Range("I58").Value = 1985.3334
$ws.Range("K58").Value = 1985.3334
$ws.Range("M58").Value = -1782.3334
$ws.Range("H86").Value = 3200
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 3200
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H136").Value = 3273
$ws.Range("I136").Value = 1985.3334
$ws.Range("K136").Value = 5956.0002
$ws.Range("M136").Value = -3406.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 58883.766
$ws.Range("J33").Value = 200039.2
$ws.Range("L33").Value = 1200235.2
$ws.Range("N33").Value = -1200801.2
$ws.Range("H55").Value = 11442.167
$ws.Range("J55").Value = 12050.546
$ws.Range("L55").Value = 36151.638
$ws.Range("N55").Value = -36505.638
$ws.Range("H129").Value = 1280
$ws.Range("J129").Value = 3250
$ws.Range("L129").Value = 9750
$ws.Range("N129").Value = -19750

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3342667
$ws.Range("I2").Value = 4007200
$ws.Range("J2").Value = 20002
$ws.Range("K2").Value = 4007200
$ws.Range("L2").Value = 20002
$ws.Range("M2").Value = -4007088
$ws.Range("N2").Value = -20226
$ws.Range("H4").Value = 3808.6667
$ws.Range("I4").Value = 3808.6667
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3808.6667
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -3695.6667
$ws.Range("N4").ClearContents()
$ws.Range("H23").Value = 48999
$ws.Range("I23").Value = 48999
$ws.Range("K23").Value = 48999
$ws.Range("M23").Value = -48769
$ws.Range("H28").Value = 3808.6667
$ws.Range("I28").Value = 3808.6667
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 3808.6667
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -3576.6667
$ws.Range("N28").ClearContents()
$ws.Range("H37").Value = 3808.6667
$ws.Range("I37").Value = 3808.6667
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 3808.6667
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -3701.6667
$ws.Range("N37").ClearContents()
$ws.Range("H46").Value = 1683.9131
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 1486.5
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 1486.5
$ws.Range("M46").Value = -2812
$ws.Range("N46").Value = -1862.5
$ws.Range("H47").Value = 24663
$ws.Range("I47").Value = 19494.5
$ws.Range("J47").Value = 35000
$ws.Range("K47").Value = 19494.5
$ws.Range("L47").Value = 35000
$ws.Range("M47").Value = -19004.5
$ws.Range("N47").Value = -35980
$ws.Range("H52").Value = 24663
$ws.Range("I52").Value = 19494.5
$ws.Range("J52").Value = 35000
$ws.Range("K52").Value = 19494.5
$ws.Range("L52").Value = 35000
$ws.Range("M52").Value = -19261.5
$ws.Range("N52").Value = -35466

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 34000
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 34000
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 10017.5
$ws.Range("J40").Value = 10017.5
$ws.Range("L40").Value = 10017.5
$ws.Range("N40").Value = -10315.5
$ws.Range("H122").Value = 4143.0586
$ws.Range("I122").Value = 4243.3
$ws.Range("J122").Value = 3999.8572
$ws.Range("K122").Value = 12729.9
$ws.Range("L122").Value = 11999.5716
$ws.Range("M122").Value = -10279.9
$ws.Range("N122").Value = -16899.5716
